$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix text labels for benchmark accuracy
$ws.Range("D1").Value = "Search Time (s)"
$ws.Range("A5").Value = "* DIAMOND"
$ws.Range("A6").Value = "* MMseqs2"

# Row 2's data cells (B2:G2) carried a stray cell-level style override;
# clear it so they fall back to their columns' default style, same as
# every other data row already does (matching columns' <col style=".."/>).
$b2 = $ws.Range("B2").Value2
$c2 = $ws.Range("C2").Value2
$d2 = $ws.Range("D2").Value2
$e2 = $ws.Range("E2").Value2
$f2 = $ws.Range("F2").Value2
$g2 = $ws.Range("G2").Value2
$ws.Range("B2:G2").ClearContents()
$ws.Range("B2").Value = $b2
$ws.Range("C2").Value = $c2
$ws.Range("D2").Value = $d2
$ws.Range("E2").Value = $e2
$ws.Range("F2").Value = $f2
$ws.Range("G2").Value = $g2

# Row heights: header row and last row shrink to 13.8
$ws.Rows.Item(1).RowHeight = 13.8
$ws.Rows.Item(6).RowHeight = 13.8

# Move selection to A6
$ws.Range("A6").Select()
